# Update LR-pair (Bmp7-Acvr2a) sheet with recomputed values based on new TPM
# expression input for the "ECs" target cluster's receptor (Acvr2a) expression.
#
# Receptor average/total expression for target cluster "ECs" roughly doubled
# (new TPM numbers), which ripples through the derived-specificity and
# edge-weight columns for every row that references that cluster's receptor
# values (rows 2 and 5), as well as the specificity columns for the other
# rows in the same "sending cluster" group (rows 3-4 and 6-7), since those
# specificities are normalized across the group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending cluster FAPs -> Target cluster ECs) ---
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 67.40485854295278
$ws.Range("R2").Value = 606.643726886575
$ws.Range("S2").Value = 0.3112096638900835
$ws.Range("T2").Value = 0.3112096638900834

# --- Row 3 (Sending cluster FAPs -> Target cluster FAPs); M/N unchanged,
#     but derived specificities shift because the ECs group total changed ---
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("S3").Value = 0.5323996960210983
$ws.Range("T3").Value = 0.532399696021098

# --- Row 4 (Sending cluster FAPs -> Target cluster MuSCs) ---
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("S4").Value = 0.1227131495228376
$ws.Range("T4").Value = 0.1227131495228375

# --- Row 5 (Sending cluster MuSCs -> Target cluster ECs) ---
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 2.349139614900555
$ws.Range("R5").Value = 21.142256534105
$ws.Range("S5").Value = 0.01084602750880658
$ws.Range("T5").Value = 0.01084602750880658

# --- Row 6 (Sending cluster MuSCs -> Target cluster FAPs) ---
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.01855476361673835
$ws.Range("T6").Value = 0.01855476361673835

# --- Row 7 (Sending cluster MuSCs -> Target cluster MuSCs) ---
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.004276699440435986
$ws.Range("T7").Value = 0.004276699440435985
